# Répartition Tâches.xlsx - "Modification de la documentation"
# Adds a new row (Intégration / Equipe Dev / 07-01-2016 / 1 day) to the
# task-tracking table, growing Table1 by one row, and moves the active
# selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Grow Table1 (A1:D15 -> A1:D16) by appending a new list row. This keeps
# the autoFilter/table ref and row striping in sync with the table.
$lo = $ws.ListObjects.Item("Table1")
$newRow = $lo.ListRows.Add()

# New task entry in row 16.
$ws.Range("A16").Value = "Intégration"
$ws.Range("B16").Value = "Equipe Dev"

# Match the date-column formatting used by the other rows (style carries
# the m/d/yyyy number format + centered alignment already defined in the
# workbook) before writing the date serial value.
$ws.Range("C7").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 42376
$ws.Application.CutCopyMode = $false

$ws.Range("D16").Value = 1

# Move the selection as recorded in the saved view state.
$ws.Range("C11").Select()
